$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("M9").Value = -15.23077000000001
$ws.Range("I9").Value = 184.23077
$ws.Range("K9").Value = 184.23077
$ws.Range("H9").Value = 189.84848
$ws.Range("I98").Value = 758.6667
$ws.Range("K98").Value = 758.6667
$ws.Range("H98").Value = 1328.8572
$ws.Range("M98").Value = 739.3333
$ws.Range("J99").Value = 200
$ws.Range("L99").Value = 600
$ws.Range("N99").Value = -3596
$ws.Range("M99").Value = 959.5
$ws.Range("H99").Value = 186.33333
$ws.Range("I99").Value = 179.5
$ws.Range("K99").Value = 538.5
$ws.Range("H112").Value = 114831.336
$ws.Range("L112").Value = 384893.25
$ws.Range("J112").Value = 128297.75
$ws.Range("N112").Value = -387109.25
$ws.Range("K116").Value = 71455470
$ws.Range("H116").Value = 62523800
$ws.Range("I116").Value = 71455470
$ws.Range("M116").Value = -71452028
$ws.Range("H122").Value = 1328.8572
$ws.Range("I122").Value = 758.6667
$ws.Range("M122").Value = 173.9998999999998
$ws.Range("K122").Value = 2276.0001
$ws.Range("M132").Value = -2029.588400000001
$ws.Range("I132").Value = 1519.8628
$ws.Range("K132").Value = 4559.588400000001
$ws.Range("H132").Value = 1491.0927
$ws.Range("N138").Value = -20269.4426
$ws.Range("H138").Value = 2859.26
$ws.Range("J138").Value = 3329.8142
$ws.Range("K138").Value = 5283.9
$ws.Range("M138").Value = -143.8999999999996
$ws.Range("I138").Value = 1761.3
$ws.Range("L138").Value = 9989.442599999998

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("K32").Value = 21447596
$ws.Range("M32").Value = -21447309
$ws.Range("I32").Value = 21447596
$ws.Range("H32").Value = 18982098
$ws.Range("J45").Value = 4528.4287
$ws.Range("N45").Value = -5282.4287
$ws.Range("L45").Value = 4528.4287
$ws.Range("K45").Value = 3351.875
$ws.Range("H45").Value = 3709.9565
$ws.Range("I45").Value = 3351.875
$ws.Range("M45").Value = -2974.875
$ws.Range("I74").Value = 2082.9302
$ws.Range("K74").Value = 2082.9302
$ws.Range("H74").Value = 2381.2156
$ws.Range("M74").Value = -1208.9302
$ws.Range("K77").Value = 10414.651
$ws.Range("H77").Value = 2381.2156
$ws.Range("M77").Value = -6046.650999999998
$ws.Range("I77").Value = 2082.9302
$ws.Range("N117").Value = -93905.5
$ws.Range("J117").Value = 84727.5
$ws.Range("L117").Value = 84727.5
$ws.Range("H117").Value = 84727.5
$ws.Range("N125").Value = -95662.2
$ws.Range("J125").Value = 85822.2
$ws.Range("H125").Value = 85822.2
$ws.Range("L125").Value = 85822.2
$ws.Range("M132").Value = -7038.9095
$ws.Range("I132").Value = 3189.6365
$ws.Range("K132").Value = 9568.9095
$ws.Range("H132").Value = 3538.1936

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("J62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("H65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("L102").Value = 90329.664
$ws.Range("N102").Value = -96819.664
$ws.Range("J102").Value = 90329.664
$ws.Range("H102").Value = 49398.168
$ws.Range("H116").Value = 87746.39999999999
$ws.Range("N116").Value = -96924.39999999999
$ws.Range("J116").Value = 87746.39999999999
$ws.Range("L116").Value = 87746.39999999999
$ws.Range("I134").Value = 7146138.5
$ws.Range("H134").Value = 4469019
$ws.Range("K134").Value = 21438415.5
$ws.Range("M134").Value = -21435880.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 758.8
$ws.Range("M22").Value = -398.5
$ws.Range("K22").Value = 748.5
$ws.Range("I22").Value = 748.5
$ws.Range("N68").Value = -66498
$ws.Range("J68").Value = 65000
$ws.Range("L68").Value = 65000
$ws.Range("H68").Value = 42666.668
$ws.Range("H71").Value = 42666.668
$ws.Range("N71").Value = -202488
$ws.Range("L71").Value = 195000
$ws.Range("J71").Value = 65000
$ws.Range("H112").Value = 60701
$ws.Range("L112").Value = 60701
$ws.Range("J112").Value = 60701
$ws.Range("N112").Value = -63655
$ws.Range("N117").Value = -57552.5
$ws.Range("J117").Value = 48374.5
$ws.Range("L117").Value = 48374.5
$ws.Range("H117").Value = 48374.5
$ws.Range("H118").Value = 115994.5
$ws.Range("L118").Value = 112989
$ws.Range("N118").Value = -116303
$ws.Range("J118").Value = 112989
$ws.Range("I122").Value = 4766719
$ws.Range("M122").Value = -14297707
$ws.Range("N122").Value = -23413.3339
$ws.Range("L122").Value = 18513.3339
$ws.Range("K122").Value = 14300157
$ws.Range("J122").Value = 6171.1113

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("M5").Value = -4552.6667
$ws.Range("K5").Value = 4664.6667
$ws.Range("I5").Value = 1554.8889
$ws.Range("H5").Value = 2135.7273
$ws.Range("K131").Value = 2426.625
$ws.Range("H131").Value = 2094.625
$ws.Range("I131").Value = 808.875
$ws.Range("M131").Value = 2613.375
$ws.Range("N132").Value = -18502.1003
$ws.Range("J132").Value = 1493.5667
$ws.Range("L132").Value = 13442.1003
$ws.Range("H132").Value = 1445.9479
$ws.Range("M135").Value = -11459.0001
$ws.Range("I135").Value = 1554.8889
$ws.Range("K135").Value = 13994.0001
$ws.Range("H135").Value = 2135.7273

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("L11").Value = 1499125
$ws.Range("M11").Value = -8003910
$ws.Range("J11").Value = 1499125
$ws.Range("K11").Value = 8004049
$ws.Range("I11").Value = 8004049
$ws.Range("N11").Value = -1499403
$ws.Range("H11").Value = 4001018.8
$ws.Range("H116").Value = 114990
$ws.Range("N116").Value = -124168
$ws.Range("J116").Value = 114990
$ws.Range("L116").Value = 114990
$ws.Range("H122").Value = 1749.6666
$ws.Range("N122").ClearContents()
$ws.Range("L122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("H124").Value = 153899
$ws.Range("L124").Value = 153899
$ws.Range("N124").Value = -163719
$ws.Range("J124").Value = 153899
$ws.Range("M126").Value = -6822.5
$ws.Range("K126").Value = 9292.5
$ws.Range("I126").Value = 3097.5
$ws.Range("H126").Value = 3231.6667
$ws.Range("H139").Value = 400000
$ws.Range("N139").Value = -410280
$ws.Range("L139").Value = 400000
$ws.Range("J139").Value = 400000

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I7").Value = 4357.706
$ws.Range("K7").Value = 4357.706
$ws.Range("H7").Value = 4613.3335
$ws.Range("M7").Value = -4245.706
$ws.Range("L7").Value = 5234.143
$ws.Range("J7").Value = 5234.143
$ws.Range("N7").Value = -5458.143
$ws.Range("J61").Value = 1198.5
$ws.Range("N61").Value = -1602.5
$ws.Range("L61").Value = 1198.5
$ws.Range("H61").Value = 1623.0834
$ws.Range("H113").Value = 1623.0834
$ws.Range("N113").Value = -5538.5
$ws.Range("J113").Value = 1198.5
$ws.Range("L113").Value = 1198.5
$ws.Range("H122").Value = 5241
$ws.Range("I122").Value = 4457.5
$ws.Range("M122").Value = -10922.5
$ws.Range("N122").Value = -31905.4
$ws.Range("L122").Value = 27005.4
$ws.Range("K122").Value = 13372.5
$ws.Range("J122").Value = 9001.799999999999
$ws.Range("M126").Value = -10603.118
$ws.Range("K126").Value = 13073.118
$ws.Range("L126").Value = 15702.429
$ws.Range("N126").Value = -20642.429
$ws.Range("I126").Value = 4357.706
$ws.Range("J126").Value = 5234.143
$ws.Range("H126").Value = 4613.3335

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("J16").Value = 111388.336
$ws.Range("H16").Value = 111388.336
$ws.Range("L16").Value = 111388.336
$ws.Range("N16").Value = -111972.336
$ws.Range("K113").Value = 2999.25
$ws.Range("H113").Value = 988.125
$ws.Range("N113").Value = -7269.5
$ws.Range("J113").Value = 976.5
$ws.Range("I113").Value = 999.75
$ws.Range("L113").Value = 2929.5
$ws.Range("M113").Value = -829.25
$ws.Range("H124").Value = 81999.7
$ws.Range("L124").Value = 81999.7
$ws.Range("N124").Value = -91819.7
$ws.Range("J124").Value = 81999.7
$ws.Range("N125").Value = -65964.875
$ws.Range("J125").Value = 56124.875
$ws.Range("H125").Value = 54666.555
$ws.Range("L125").Value = 56124.875
$ws.Range("J135").Value = 10000
$ws.Range("L135").Value = 10000
$ws.Range("N135").Value = -20140
$ws.Range("H135").Value = 10000
